$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 11.6312914161566
$ws.Cells.Item(2, 3).Value = 8.015445958818585
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 12.25441670958022
$ws.Cells.Item(2, 6).Value = 16.86991607391233
$ws.Cells.Item(2, 7).Value = 3.662559065045042
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 25.03053132189126
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 9.029295991837998
$ws.Cells.Item(2, 12).Value = 10.15121155305382
$ws.Cells.Item(2, 13).Value = 13.93153174997595
$ws.Cells.Item(2, 14).Value = 20.1223370094384
$ws.Cells.Item(2, 15).Value = 24.96835095978425

$ws.Cells.Item(3, 2).Value = 11.42600299865751
$ws.Cells.Item(3, 3).Value = 7.990342930504857
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = 12.2792371126617
$ws.Cells.Item(3, 6).Value = 15.89584955866808
$ws.Cells.Item(3, 7).Value = 3.664058736067018
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 25.12550251550319
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 8.884426446292636
$ws.Cells.Item(3, 12).Value = 10.15873447998605
$ws.Cells.Item(3, 13).Value = 13.90287015752963
$ws.Cells.Item(3, 14).Value = 20.17880671878433
$ws.Cells.Item(3, 15).Value = 25.05905202467667

$ws.Cells.Item(4, 2).Value = 11.29995938865098
$ws.Cells.Item(4, 3).Value = 7.974746885694811
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 12.29609681548011
$ws.Cells.Item(4, 6).Value = 15.26997757108491
$ws.Cells.Item(4, 7).Value = 3.665028865048658
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 25.1877164681142
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 8.795551901708979
$ws.Cells.Item(4, 12).Value = 10.16470999796072
$ws.Cells.Item(4, 13).Value = 13.88711623537405
$ws.Cells.Item(4, 14).Value = 20.21517809924031
$ws.Cells.Item(4, 15).Value = 25.11906344312224

$ws.Cells.Item(5, 2).Value = 11.24866381075409
$ws.Cells.Item(5, 3).Value = 7.968346790914143
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = 12.3033749628227
$ws.Cells.Item(5, 6).Value = 15.00819731993403
$ws.Cells.Item(5, 7).Value = 3.665436640690587
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 25.21405075309706
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 8.759397843310014
$ws.Cells.Item(5, 12).Value = 10.16748677419599
$ws.Cells.Item(5, 13).Value = 13.88116492207993
$ws.Cells.Item(5, 14).Value = 20.23042807392745
$ws.Cells.Item(5, 15).Value = 25.14460470130836

$ws.Cells.Item(6, 2).Value = 11.24015231938786
$ws.Cells.Item(6, 3).Value = 7.967281406971008
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = 12.30460812648951
$ws.Cells.Item(6, 6).Value = 14.96433081551593
$ws.Cells.Item(6, 7).Value = 3.665505103973679
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = 25.21848284481609
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 8.75339958889375
$ws.Cells.Item(6, 12).Value = 10.16796851078672
$ws.Cells.Item(6, 13).Value = 13.88020514213158
$ws.Cells.Item(6, 14).Value = 20.23298622608302
$ws.Cells.Item(6, 15).Value = 25.14891139715697

$ws.Cells.Item(7, 2).Value = 11.29926723294973
$ws.Cells.Item(7, 3).Value = 7.974660750492256
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(7, 5).Value = 12.2961933199761
$ws.Cells.Item(7, 6).Value = 15.26647399323137
$ws.Cells.Item(7, 7).Value = 3.665034314033839
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 25.1880676462482
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 8.795064001834499
$ws.Cells.Item(7, 12).Value = 10.16474606221058
$ws.Cells.Item(7, 13).Value = 13.88703407075094
$ws.Cells.Item(7, 14).Value = 20.21538202982717
$ws.Cells.Item(7, 15).Value = 25.11940350424281

$ws.Cells.Item(8, 2).Value = 11.56054702365648
$ws.Cells.Item(8, 3).Value = 8.006829164592451
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 12.26263879612153
$ws.Cells.Item(8, 6).Value = 16.53996406344768
$ws.Cells.Item(8, 7).Value = 3.663065936572361
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).Value = 25.06246803816506
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 8.979355426886562
$ws.Cells.Item(8, 12).Value = 10.15352432763163
$ws.Cells.Item(8, 13).Value = 13.9212693613469
$ws.Cells.Item(8, 14).Value = 20.14145590059941
$ws.Cells.Item(8, 15).Value = 24.99872774912244

$ws.Cells.Item(9, 2).Value = 12.07006310102962
$ws.Cells.Item(9, 3).Value = 8.068403370964058
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 12.20967556638052
$ws.Cells.Item(9, 6).Value = 19.00274580682531
$ws.Cells.Item(9, 7).Value = 3.659595632788434
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 24.84709542968811
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 9.339452111720542
$ws.Cells.Item(9, 12).Value = 10.14225182333835
$ws.Cells.Item(9, 13).Value = 14.00281889802661
$ws.Cells.Item(9, 14).Value = 20.00991138862782
$ws.Cells.Item(9, 15).Value = 24.79637961075894

$ws.Cells.Item(10, 2).Value = 12.43892074167141
$ws.Cells.Item(10, 3).Value = 8.112646062053379
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = 12.17856718984426
$ws.Cells.Item(10, 6).Value = 20.67494806633232
$ws.Cells.Item(10, 7).Value = 3.65728118739446
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 24.70767822386784
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 9.600741110532459
$ws.Cells.Item(10, 12).Value = 10.14047121232617
$ws.Cells.Item(10, 13).Value = 14.07120837622394
$ws.Cells.Item(10, 14).Value = 19.92137233539386
$ws.Cells.Item(10, 15).Value = 24.66864234003482

$ws.Cells.Item(11, 2).Value = 12.60478424717608
$ws.Cells.Item(11, 3).Value = 8.132540563791201
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(11, 5).Value = 12.16610454617019
$ws.Cells.Item(11, 6).Value = 21.3917225636224
$ws.Cells.Item(11, 7).Value = 3.656278850777209
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 24.64833229645113
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 9.7183956086204
$ws.Cells.Item(11, 12).Value = 10.14106248749925
$ws.Cells.Item(11, 13).Value = 14.10408822999062
$ws.Cells.Item(11, 14).Value = 19.88283750719711
$ws.Cells.Item(11, 15).Value = 24.61508015332619

$ws.Cells.Item(12, 2).Value = 12.66725587750028
$ws.Cells.Item(12, 3).Value = 8.140039312656334
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(12, 5).Value = 12.16162763851729
$ws.Cells.Item(12, 6).Value = 21.65686569030329
$ws.Cells.Item(12, 7).Value = 3.655906518123708
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 9).Value = 24.62644519300849
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 9.762734915301897
$ws.Cells.Item(12, 12).Value = 10.14148686827105
$ws.Cells.Item(12, 13).Value = 14.1167867495
$ws.Cells.Item(12, 14).Value = 19.8684946641152
$ws.Cells.Item(12, 15).Value = 24.59545175575104

$ws.Cells.Item(13, 2).Value = 12.6538174097998
$ws.Cells.Item(13, 3).Value = 8.138425898332914
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 12.16258104559555
$ws.Cells.Item(13, 6).Value = 21.60004134736742
$ws.Cells.Item(13, 7).Value = 3.655986385578444
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 24.63113291606007
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 9.753195788169807
$ws.Cells.Item(13, 12).Value = 10.14138657114834
$ws.Cells.Item(13, 13).Value = 14.1140409931882
$ws.Cells.Item(13, 14).Value = 19.8715725736598
$ws.Cells.Item(13, 15).Value = 24.59964996460417

$ws.Cells.Item(14, 2).Value = 12.60993091540746
$ws.Cells.Item(14, 3).Value = 8.133158198042958
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 12.16573137217274
$ws.Cells.Item(14, 6).Value = 21.4136618050453
$ws.Cells.Item(14, 7).Value = 3.656248074027932
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 9).Value = 24.64651988957015
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 9.722047941419319
$ws.Cells.Item(14, 12).Value = 10.14109338984264
$ws.Cells.Item(14, 13).Value = 14.10512802162114
$ws.Cells.Item(14, 14).Value = 19.88165251923902
$ws.Cells.Item(14, 15).Value = 24.61345219320325

$ws.Cells.Item(15, 2).Value = 12.58300353243268
$ws.Cells.Item(15, 3).Value = 8.129926992779005
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 12.16769259641426
$ws.Cells.Item(15, 6).Value = 21.29868154950795
$ws.Cells.Item(15, 7).Value = 3.656409306385429
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 24.65602115375571
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 9.702939964231112
$ws.Cells.Item(15, 12).Value = 10.14093988416358
$ws.Cells.Item(15, 13).Value = 14.09970061282218
$ws.Cells.Item(15, 14).Value = 19.88785923330357
$ws.Cells.Item(15, 15).Value = 24.62199170534107

$ws.Cells.Item(16, 2).Value = 12.42803714559413
$ws.Cells.Item(16, 3).Value = 8.111341114767859
$ws.Cells.Item(16, 4).Value = 0
$ws.Cells.Item(16, 5).Value = 12.17941559975139
$ws.Cells.Item(16, 6).Value = 20.62722412089977
$ws.Cells.Item(16, 7).Value = 3.65734770607766
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 9).Value = 24.71163861778836
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 9.593024334939432
$ws.Cells.Item(16, 12).Value = 10.1404606826875
$ws.Cells.Item(16, 13).Value = 14.06909462770226
$ws.Cells.Item(16, 14).Value = 19.9239256418195
$ws.Cells.Item(16, 15).Value = 24.67223429549637

$ws.Cells.Item(17, 2).Value = 12.33243224790733
$ws.Cells.Item(17, 3).Value = 8.099878833695833
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = 12.18703951980774
$ws.Cells.Item(17, 6).Value = 20.20408069597325
$ws.Cells.Item(17, 7).Value = 3.657936298233843
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = 24.74680187839797
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 11).Value = 9.525256234714998
$ws.Cells.Item(17, 12).Value = 10.14052490621377
$ws.Cells.Item(17, 13).Value = 14.05076723710988
$ws.Cells.Item(17, 14).Value = 19.94649664603662
$ws.Cells.Item(17, 15).Value = 24.70422137662845

$ws.Cells.Item(18, 2).Value = 12.27726439427942
$ws.Cells.Item(18, 3).Value = 8.093264315157024
$ws.Cells.Item(18, 4).Value = 0
$ws.Cells.Item(18, 5).Value = 12.19158357291822
$ws.Cells.Item(18, 6).Value = 19.95656407809801
$ws.Cells.Item(18, 7).Value = 3.658279597694512
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 24.76741048397085
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 11).Value = 9.486166494193386
$ws.Cells.Item(18, 12).Value = 10.14069380771897
$ws.Cells.Item(18, 13).Value = 14.04039274887946
$ws.Cells.Item(18, 14).Value = 19.95964292041378
$ws.Cells.Item(18, 15).Value = 24.72304736087272

$ws.Cells.Item(19, 2).Value = 12.25855672114476
$ws.Cells.Item(19, 3).Value = 8.091021069795913
$ws.Cells.Item(19, 4).Value = 0
$ws.Cells.Item(19, 5).Value = 12.19314942751505
$ws.Cells.Item(19, 6).Value = 19.87204792380568
$ws.Cells.Item(19, 7).Value = 3.658396651012558
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 9).Value = 24.7744540977053
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 9.472913517084677
$ws.Cells.Item(19, 12).Value = 10.1407736893794
$ws.Cells.Item(19, 13).Value = 14.03690899964139
$ws.Cells.Item(19, 14).Value = 19.96412222725071
$ws.Cells.Item(19, 15).Value = 24.72949497966387

$ws.Cells.Item(20, 2).Value = 12.34262846150124
$ws.Cells.Item(20, 3).Value = 8.101101271958811
$ws.Cells.Item(20, 4).Value = 0
$ws.Cells.Item(20, 5).Value = 12.18621148997424
$ws.Cells.Item(20, 6).Value = 20.24955283636154
$ws.Cells.Item(20, 7).Value = 3.657873149525975
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 9).Value = 24.7430189864875
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 9.532482083645959
$ws.Cells.Item(20, 12).Value = 10.1405044191339
$ws.Cells.Item(20, 13).Value = 14.05270098941026
$ws.Cells.Item(20, 14).Value = 19.94407695522642
$ws.Cells.Item(20, 15).Value = 24.70077201006694

$ws.Cells.Item(21, 2).Value = 12.6228310554432
$ws.Cells.Item(21, 3).Value = 8.134706407949029
$ws.Cells.Item(21, 4).Value = 0
$ws.Cells.Item(21, 5).Value = 12.1647994689539
$ws.Cells.Item(21, 6).Value = 21.46857628470577
$ws.Cells.Item(21, 7).Value = 3.656171013821404
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 9).Value = 24.64198446215046
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 11).Value = 9.731202936536201
$ws.Cells.Item(21, 12).Value = 10.14117407195535
$ws.Cells.Item(21, 13).Value = 14.10773931200623
$ws.Cells.Item(21, 14).Value = 19.87868503285058
$ws.Cells.Item(21, 15).Value = 24.60938037709015

$ws.Cells.Item(22, 2).Value = 12.80396830160978
$ws.Cells.Item(22, 3).Value = 8.15646535285253
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(22, 5).Value = 12.15221830365885
$ws.Cells.Item(22, 6).Value = 22.22866616901552
$ws.Cells.Item(22, 7).Value = 3.65510069902601
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 9).Value = 24.57936764529455
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 9.859814329611732
$ws.Cells.Item(22, 12).Value = 10.14277973465846
$ws.Cells.Item(22, 13).Value = 14.14515010463411
$ws.Cells.Item(22, 14).Value = 19.83740120052517
$ws.Cells.Item(22, 15).Value = 24.55346557201587

$ws.Cells.Item(23, 2).Value = 12.7074930183106
$ws.Cells.Item(23, 3).Value = 8.144871368660359
$ws.Cells.Item(23, 4).Value = 0
$ws.Cells.Item(23, 5).Value = 12.15880397724508
$ws.Cells.Item(23, 6).Value = 21.82633154458858
$ws.Cells.Item(23, 7).Value = 3.655668102400526
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 9).Value = 24.61247498437803
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 11).Value = 9.791300528717874
$ws.Cells.Item(23, 12).Value = 10.14181624881297
$ws.Cells.Item(23, 13).Value = 14.12505379825324
$ws.Cells.Item(23, 14).Value = 19.8593024982429
$ws.Cells.Item(23, 15).Value = 24.58295905532859

$ws.Cells.Item(24, 2).Value = 12.33801938208206
$ws.Cells.Item(24, 3).Value = 8.100548684451773
$ws.Cells.Item(24, 4).Value = 0
$ws.Cells.Item(24, 5).Value = 12.18658534086707
$ws.Cells.Item(24, 6).Value = 20.22900810905287
$ws.Cells.Item(24, 7).Value = 3.65790168376034
$ws.Cells.Item(24, 8).Value = 0
$ws.Cells.Item(24, 9).Value = 24.74472800842973
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 11).Value = 9.529215675925252
$ws.Cells.Item(24, 12).Value = 10.14051327017091
$ws.Cells.Item(24, 13).Value = 14.05182623433508
$ws.Cells.Item(24, 14).Value = 19.94517036812097
$ws.Cells.Item(24, 15).Value = 24.70233011005538

$ws.Cells.Item(25, 2).Value = 11.93292576291023
$ws.Cells.Item(25, 3).Value = 8.051915212651375
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 12.22263141325678
$ws.Cells.Item(25, 6).Value = 18.34778573295695
$ws.Cells.Item(25, 7).Value = 3.660492970401134
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).Value = 24.90205237216912
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 9.242436167088718
$ws.Cells.Item(25, 12).Value = 10.14415630281993
$ws.Cells.Item(25, 13).Value = 13.97924547092837
$ws.Cells.Item(25, 14).Value = 20.04406851219702
$ws.Cells.Item(25, 15).Value = 24.84744644792568

